# Update crypto price/volume data per the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.639.19"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = "'3.232.62"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('D5').Value = "'605.91"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('D6').Value = "'157.95"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.30%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = "'3.230.88"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.74%  '
$ws.Range('D9').Value = "'0.549"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('D10').Value = "'0.162"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.00%  '
$ws.Range('D11').Value = "'5.70"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.16%  '
$ws.Range('D12').Value = "'0.507"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').Value = "'0.0000275"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.75%  '
$ws.Range('D14').Value = "'39.02"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('D15').Value = "'3.760.68"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').Value = "'66.669.12"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('D17').Value = "'7.40"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = "'3.234.10"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('D20').Value = "'511.00"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').Value = "'0.735"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('D23').Value = "'8.06"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('D24').Value = "'14.67"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('D25').Value = "'84.93"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('D27').Value = "'3.00"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('E29').Value = '  +5.17%  '
$ws.Range('E30').Value = '  +2.93%  '
$ws.Range('D31').Value = "'7.03"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('D32').Value = "'28.25"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  -2.95%  '
$ws.Range('D35').Value = "'0.104"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +16.89%  '
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').Value = "'509.43"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.62%  '
$ws.Range('D38').Value = "'55.79"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.06%  '
$ws.Range('D39').Value = "'0.0₃0776"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +18.98%  '
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('D41').Value = "'3.07"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.54%  '
$ws.Range('E42').Value = '  +6.17%  '
$ws.Range('D43').Value = "'8.75"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.15%  '
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('E45').Value = '  +2.50%  '
$ws.Range('D46').Value = "'2.870.44"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').Value = "'28.55"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('E48').Value = '  +4.85%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').Value = "'122.40"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.57%  '
